# edit.ps1 - apply the changes described by the diff:
#  1. Change the table's style (tableStyleId) on slide 5 from
#     {7CCA818E-CA92-45FF-9684-036A0058463B} to {47E78101-1DA3-43A4-840E-ABAD3CD8219A}
#  2. Swap the "Integral" (Red Violet) theme colors for the "Office Theme" colors
#     on the deck's master theme (ppt/theme/theme1.xml) - the font scheme and the
#     format scheme are already identical between the two themes, only the color
#     scheme differs, so updating the 12 theme colors reproduces that half of the
#     swap.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{47E78101-1DA3-43A4-840E-ABAD3CD8219A}")

# --- 2. Theme colors --------------------------------------------------------
function ConvertTo-RGBInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

# Office theme color order matches ThemeColorScheme.Colors(1..12):
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-RGBInt($officeColors[$i - 1])
}
